# Update the "Förändrad" (Changed) date column (C) for rows 2-33
# from 45221 (2023-10-22) to 45224 (2023-10-25), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
